# Insert the HP ProLiant (Xeon E5-2666G) benchmark numbers that were
# previously left as placeholder zeros in row 4 (columns B:G).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B4").Value = 0.01
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 0.4
$ws.Range("F4").Value = 5.31
$ws.Range("G4").Value = 53

# Leave the selection where the author ended up after entering the data.
$null = $ws.Range("G8").Select()
